$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.3069574617629485 ; $ws.Range("C2").Value = 0.04707455330151333 ; $ws.Range("D2").Value = 0.03034126822113592 ; $ws.Range("E2").Value = 0.1597681957234656 ; $ws.Range("F2").Value = 0.746648779129977 ; $ws.Range("H2").Value = 0.07973214163530429 ; $ws.Range("K2").Value = 0.2775569654752985 ; $ws.Range("M2").Value = 0.2255241579484988 ; $ws.Range("N2").Value = 1.50133500787515 ; $ws.Range("O2").Value = 2.576247709650886
$ws.Range("B3").Value = 0.2734750974589133 ; $ws.Range("C3").Value = 0.04344008859813187 ; $ws.Range("D3").Value = 0.02844548043131567 ; $ws.Range("E3").Value = 0.1487329861241093 ; $ws.Range("F3").Value = 0.7437578374085732 ; $ws.Range("H3").Value = 0.07973214163530429 ; $ws.Range("K3").Value = 0.2431398117498276 ; $ws.Range("M3").Value = 0.2033048863317291 ; $ws.Range("N3").Value = 1.518021207761477 ; $ws.Range("O3").Value = 2.579888993381758
$ws.Range("B4").Value = 0.2529545563172064 ; $ws.Range("C4").Value = 0.04119035834405338 ; $ws.Range("D4").Value = 0.02726995512390573 ; $ws.Range("E4").Value = 0.142063879495808 ; $ws.Range("F4").Value = 0.7424312632611816 ; $ws.Range("H4").Value = 0.07973214163530429 ; $ws.Range("K4").Value = 0.2219865104311367 ; $ws.Range("M4").Value = 0.1897423394864575 ; $ws.Range("N4").Value = 1.528791050497812 ; $ws.Range("O4").Value = 2.583585944296601
$ws.Range("B5").Value = 0.2446021479261162 ; $ws.Range("C5").Value = 0.040269044277359 ; $ws.Range("D5").Value = 0.02678805104181947 ; $ws.Range("E5").Value = 0.139372760029211 ; $ws.Range("F5").Value = 0.7420034464434551 ; $ws.Range("H5").Value = 0.07973214163530429 ; $ws.Range("K5").Value = 0.2133614458602153 ; $ws.Range("M5").Value = 0.1842356117007142 ; $ws.Range("N5").Value = 1.533311674359589 ; $ws.Range("O5").Value = 2.58545981243671
$ws.Range("B6").Value = 0.2432158447921609 ; $ws.Range("C6").Value = 0.04011578814649397 ; $ws.Range("D6").Value = 0.02670785874741455 ; $ws.Range("E6").Value = 0.1389275021948464 ; $ws.Range("F6").Value = 0.7419392190108738 ; $ws.Range("H6").Value = 0.07973214163530429 ; $ws.Range("K6").Value = 0.2119289739284369 ; $ws.Range("M6").Value = 0.1833224377961287 ; $ws.Range("N6").Value = 1.534070280978523 ; $ws.Range("O6").Value = 2.58579314989467
$ws.Range("B7").Value = 0.2528418723109667 ; $ws.Range("C7").Value = 0.04117795148087566 ; $ws.Range("D7").Value = 0.02726346757599174 ; $ws.Range("E7").Value = 0.1420274786954607 ; $ws.Range("F7").Value = 0.7424250369616701 ; $ws.Range("H7").Value = 0.07973214163530429 ; $ws.Range("K7").Value = 0.2218702092406346 ; $ws.Range("M7").Value = 0.1896679924168936 ; $ws.Range("N7").Value = 1.52885148352503 ; $ws.Range("O7").Value = 2.583609728762625
$ws.Range("B8").Value = 0.2954051231704682 ; $ws.Range("C8").Value = 0.04582517733281577 ; $ws.Range("D8").Value = 0.02969000258322296 ; $ws.Range("E8").Value = 0.1559410164013286 ; $ws.Range("F8").Value = 0.7455588981920656 ; $ws.Range("H8").Value = 0.07973214163530429 ; $ws.Range("K8").Value = 0.2656945045560803 ; $ws.Range("M8").Value = 0.2178462812302584 ; $ws.Range("N8").Value = 1.506979530991156 ; $ws.Range("O8").Value = 2.577199889859628
$ws.Range("B9").Value = 0.3791585438698633 ; $ws.Range("C9").Value = 0.05479322258462105 ; $ws.Range("D9").Value = 0.03435626375912904 ; $ws.Range("E9").Value = 0.1840813562574581 ; $ws.Range("F9").Value = 0.7552642961534701 ; $ws.Range("H9").Value = 0.07973214163530429 ; $ws.Range("K9").Value = 0.3514550262920295 ; $ws.Range("M9").Value = 0.2737447342510393 ; $ws.Range("N9").Value = 1.468251585543369 ; $ws.Range("O9").Value = 2.576231566654855
$ws.Range("B10").Value = 0.4408562900358675 ; $ws.Range("C10").Value = 0.0612926668104592 ; $ws.Range("D10").Value = 0.03772748102135637 ; $ws.Range("E10").Value = 0.2052944720597338 ; $ws.Range("F10").Value = 0.7645691733932551 ; $ws.Range("H10").Value = 0.07973214163530429 ; $ws.Range("K10").Value = 0.414344360189034 ; $ws.Range("M10").Value = 0.3152148378202781 ; $ws.Range("N10").Value = 1.442337222943431 ; $ws.Range("O10").Value = 2.582607475424084
$ws.Range("B11").Value = 0.4689578846499955 ; $ws.Range("C11").Value = 0.06422987158218518 ; $ws.Range("D11").Value = 0.0392485661843196 ; $ws.Range("E11").Value = 0.215065472280763 ; $ws.Range("F11").Value = 0.7692753502322063 ; $ws.Range("H11").Value = 0.07973214163530429 ; $ws.Range("K11").Value = 0.4429267800221908 ; $ws.Range("M11").Value = 0.3341704033445012 ; $ws.Range("N11").Value = 1.431099715214284 ; $ws.Range("O11").Value = 2.587050496684242
$ws.Range("B12").Value = 0.479603950741506 ; $ws.Range("C12").Value = 0.06533929441859243 ; $ws.Range("D12").Value = 0.03982274263206165 ; $ws.Range("E12").Value = 0.2187831522937529 ; $ws.Range("F12").Value = 0.7711255631686527 ; $ws.Range("H12").Value = 0.07973214163530429 ; $ws.Range("K12").Value = 0.4537461332288331 ; $ws.Range("M12").Value = 0.3413615212004117 ; $ws.Range("N12").Value = 1.426923668676897 ; $ws.Range("O12").Value = 2.588954982172083
$ws.Range("B13").Value = 0.477310932067752 ; $ws.Range("C13").Value = 0.06510048698270055 ; $ws.Range("D13").Value = 0.039699165099222 ; $ws.Range("E13").Value = 0.2179816966768726 ; $ws.Range("F13").Value = 0.7707240582585371 ; $ws.Range("H13").Value = 0.07973214163530429 ; $ws.Range("K13").Value = 0.4514161848227047 ; $ws.Range("M13").Value = 0.3398122041824578 ; $ws.Range("N13").Value = 1.427819523939984 ; $ws.Range("O13").Value = 2.588534940221109
$ws.Range("B14").Value = 0.4698336531347991 ; $ws.Range("C14").Value = 0.06432120146727982 ; $ws.Range("D14").Value = 0.03929584077383907 ; $ws.Range("E14").Value = 0.2153709738462055 ; $ws.Range("F14").Value = 0.7694262037092443 ; $ws.Range("H14").Value = 0.07973214163530429 ; $ws.Range("K14").Value = 0.4438169810437387 ; $ws.Range("M14").Value = 0.3347617591111245 ; $ws.Range("N14").Value = 1.430754558338904 ; $ws.Range("O14").Value = 2.587202729373303
$ws.Range("B15").Value = 0.4652541905030887 ; $ws.Range("C15").Value = 0.06384349649307808 ; $ws.Range("D15").Value = 0.03904855440777055 ; $ws.Range("E15").Value = 0.2137741307959331 ; $ws.Range("F15").Value = 0.7686400976983236 ; $ws.Range("H15").Value = 0.07973214163530429 ; $ws.Range("K15").Value = 0.4391616915106056 ; $ws.Range("M15").Value = 0.3316699168271384 ; $ws.Range("N15").Value = 1.432562689200186 ; $ws.Range("O15").Value = 2.586415629542103
$ws.Range("B16").Value = 0.4390204520542795 ; $ws.Range("C16").Value = 0.06110031930379023 ; $ws.Range("D16").Value = 0.03762782061459546 ; $ws.Range("E16").Value = 0.2046583677667186 ; $ws.Range("F16").Value = 0.7642711402333902 ; $ws.Range("H16").Value = 0.07973214163530429 ; $ws.Range("K16").Value = 0.4124758681930984 ; $ws.Range("M16").Value = 0.3139778726423543 ; $ws.Range("N16").Value = 1.443082718797259 ; $ws.Range("O16").Value = 2.582348170388258
$ws.Range("B17").Value = 0.422935558142882 ; $ws.Range("C17").Value = 0.05941246556781721 ; $ws.Range("D17").Value = 0.03675302520877466 ; $ws.Range("E17").Value = 0.1990972934537893 ; $ws.Range("F17").Value = 0.7617121748877622 ; $ws.Range("H17").Value = 0.07973214163530429 ; $ws.Range("K17").Value = 0.3960979615613951 ; $ws.Range("M17").Value = 0.3031476030183953 ; $ws.Range("N17").Value = 1.44967763565079 ; $ws.Range("O17").Value = 2.580248143317021
$ws.Range("B18").Value = 0.4136872707341581 ; $ws.Range("C18").Value = 0.05843983072823278 ; $ws.Range("D18").Value = 0.03624869130015895 ; $ws.Range("E18").Value = 0.1959100917405294 ; $ws.Range("F18").Value = 0.760284877692456 ; $ws.Range("H18").Value = 0.07973214163530429 ; $ws.Range("K18").Value = 0.3866753764797579 ; $ws.Range("M18").Value = 0.296926845199124 ; $ws.Range("N18").Value = 1.453522719494304 ; $ws.Range("O18").Value = 2.57918544780847
$ws.Range("B19").Value = 0.4105565424543443 ; $ws.Range("C19").Value = 0.05811020087199381 ; $ws.Range("D19").Value = 0.03607773160423733 ; $ws.Range("E19").Value = 0.1948329083922573 ; $ws.Range("F19").Value = 0.7598092706155128 ; $ws.Range("H19").Value = 0.07973214163530429 ; $ws.Range("K19").Value = 0.3834846435055113 ; $ws.Range("M19").Value = 0.2948220671307382 ; $ws.Range("N19").Value = 1.454833502841264 ; $ws.Range("O19").Value = 2.578850567240693
$ws.Range("B20").Value = 0.4246474822422783 ; $ws.Range("C20").Value = 0.0595923297515526 ; $ws.Range("D20").Value = 0.03684627045743838 ; $ws.Range("E20").Value = 0.199688099897763 ; $ws.Range("F20").Value = 0.761979970005747 ; $ws.Range("H20").Value = 0.07973214163530429 ; $ws.Range("K20").Value = 0.3978416738671342 ; $ws.Range("M20").Value = 0.3042996210159288 ; $ws.Range("N20").Value = 1.448970227160634 ; $ws.Range("O20").Value = 2.580456667367258
$ws.Range("B21").Value = 0.4720297897238197 ; $ws.Range("C21").Value = 0.06455017371231975 ; $ws.Range("D21").Value = 0.03941435667330495 ; $ws.Range("E21").Value = 0.2161373267239313 ; $ws.Range("F21").Value = 0.7698055670631447 ; $ws.Range("H21").Value = 0.07973214163530429 ; $ws.Range("K21").Value = 0.446049168383496 ; $ws.Range("M21").Value = 0.3362448426873001 ; $ws.Range("N21").Value = 1.429890313052701 ; $ws.Range("O21").Value = 2.587588005219203
$ws.Range("B22").Value = 0.503023492289941 ; $ws.Range("C22").Value = 0.06777388091623493 ; $ws.Range("D22").Value = 0.04108209599902324 ; $ws.Range("E22").Value = 0.2269906132914201 ; $ws.Range("F22").Value = 0.7753169126578996 ; $ws.Range("H22").Value = 0.07973214163530429 ; $ws.Range("K22").Value = 0.4775309447405505 ; $ws.Range("M22").Value = 0.357199007631813 ; $ws.Range("N22").Value = 1.417883074047758 ; $ws.Range("O22").Value = 2.593542876277581
$ws.Range("B23").Value = 0.4864792752467224 ; $ws.Range("C23").Value = 0.06605485330644001 ; $ws.Range("D23").Value = 0.04019297625616503 ; $ws.Range("E23").Value = 0.221188538254566 ; $ws.Range("F23").Value = 0.7723390841501327 ; $ws.Range("H23").Value = 0.07973214163530429 ; $ws.Range("K23").Value = 0.4607309166046889 ; $ws.Range("M23").Value = 0.3460084019974445 ; $ws.Range("N23").Value = 1.424249199263498 ; $ws.Range("O23").Value = 2.590246174757709
$ws.Range("B24").Value = 0.4238735240282665 ; $ws.Range("C24").Value = 0.0595110202076512 ; $ws.Range("D24").Value = 0.03680411864729649 ; $ws.Range("E24").Value = 0.1994209653793106 ; $ws.Range("F24").Value = 0.7618587631314142 ; $ws.Range("H24").Value = 0.07973214163530429 ; $ws.Range("K24").Value = 0.3970533623996459 ; $ws.Range("M24").Value = 0.3037787759032895 ; $ws.Range("N24").Value = 1.449289879649225 ; $ws.Range("O24").Value = 2.580361943114212
$ws.Range("B25").Value = 0.3564714194400835 ; $ws.Range("C25").Value = 0.05238276043012036 ; $ws.Range("D25").Value = 0.03310387632515699 ; $ws.Range("E25").Value = 0.1763752981699085 ; $ws.Range("F25").Value = 0.752257190075035 ; $ws.Range("H25").Value = 0.07973214163530429 ; $ws.Range("K25").Value = 0.3282746989990812 ; $ws.Range("M25").Value = 0.2585529892774687 ; $ws.Range("N25").Value = 1.427819523939984 ; $ws.Range("O25").Value = 2.588534940221109
